$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 134
$ws1.Range("F3").Value = 332
$ws1.Range("F5").Value = 1717
$ws1.Range("F6").Value = 80
$ws1.Range("F7").Value = 2170
$ws1.Range("F11").Value = 4866
$ws1.Range("F14").Value = 302
$ws1.Range("F15").Value = 225
$ws1.Range("F16").Value = 30
$ws1.Range("F17").Value = 174
$ws1.Range("F20").Value = 120
$ws1.Range("F21").Value = 3800
$ws1.Range("F22").Value = 698
$ws1.Range("F23").Value = 644
$ws1.Range("F26").Value = 102
$ws1.Range("F28").Value = 20
$ws1.Range("F30").Value = 84
$ws1.Range("F34").Value = 913
$ws1.Range("F35").Value = 2429

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 134
$ws4.Range("F3").Value = 332
$ws4.Range("F5").Value = 1717
$ws4.Range("F6").Value = 80
$ws4.Range("F7").Value = 2170
$ws4.Range("F11").Value = 4866
$ws4.Range("F14").Value = 302
$ws4.Range("F15").Value = 225
$ws4.Range("F16").Value = 30
$ws4.Range("F17").Value = 174
$ws4.Range("F20").Value = 120
$ws4.Range("F21").Value = 3800
$ws4.Range("F22").Value = 698
$ws4.Range("F23").Value = 644
$ws4.Range("F26").Value = 102
$ws4.Range("F28").Value = 20
$ws4.Range("F30").Value = 84
$ws4.Range("F35").Value = 913
$ws4.Range("F36").Value = 2429
